$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the pre-filled "0" value in B2 (Absent count for bunk B1)
$ws.Range("B2").ClearContents()

# Move the active selection to B2
$ws.Range("B2").Select()
